{"js": "// Apply the Review_221.docx edit:\n//  P1: date 12.06.24 -> 11.06.24\n//  P2: title replaced (new paper title)\n//  P3: body paragraph replaced\n//  P4: body paragraph replaced\n//  P5: body paragraph replaced\n//  P6: body paragraph removed entirely (old \"\u05d1\u05d3\u05f4\u05db \u05d4\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1...\" paragraph)\n//  P7 (becomes P6 after removal): URL replaced with the new arxiv link\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) Date line\nitems[0].insertText(\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 11.06.24:\u26a1\ufe0f\ud83d\ude80\",\n  \"Replace\"\n);\n\n// 2) Paper title\nitems[1].insertText(\n  \"The Geometry of Categorical and Hierarchical Concepts in Large Language Models\",\n  \"Replace\"\n);\n\n// 3) First body paragraph\nitems[2].insertText(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 \u05db\u05d9\u05e6\u05d3 \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05d5\u05de\u05d5\u05e9\u05d2\u05d9\u05dd \u05de\u05e7\u05d5\u05d3\u05d3\u05d9\u05dd \u05d1\u05de\u05e8\u05d7\u05d1\u05d9 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 (embeddings)\u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05dc \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d4. \u05d4\u05db\u05d5\u05ea\u05d1\u05d9\u05dd \u05d7\u05d5\u05e7\u05e8\u05d9\u05dd 2 \u05e9\u05d0\u05dc\u05d5\u05ea \u05de\u05e8\u05db\u05d6\u05d9\u05d5\u05ea: \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d9\u05dd \u05d5\u05d4\u05e7\u05d9\u05d3\u05d5\u05d3 \u05e9\u05dc \u05d9\u05d7\u05e1\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05d1\u05d9\u05df \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd.\",\n  \"Replace\"\n);\n\n// 4) Second body paragraph\nitems[3].insertText(\n  \"\u05d4\u05dd \u05de\u05e8\u05d7\u05d9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05d4\u05e1\u05ea\u05db\u05dc\u05d5\u05ea \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d4\u05e8\u05d2\u05d9\u05dc\u05d4 \u05e2\u05dc \u05d4\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05d4\u05e8\u05d0\u05d5\u05ea \u05e9\u05d4\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d9\u05dd \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05db\u05e1\u05d9\u05de\u05e4\u05dc\u05e7\u05e1\u05d9\u05dd, \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05d4\u05dd \u05d0\u05d5\u05e8\u05ea\u05d5\u05d2\u05d5\u05e0\u05dc\u05d9\u05d9\u05dd, \u05d5\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05db\u05e4\u05d5\u05dc\u05d9\u05d8\u05d5\u05e4\u05d9\u05dd \u05e9\u05e0\u05d1\u05e0\u05d9\u05dd \u05de\u05e1\u05db\u05d5\u05de\u05d9\u05dd \u05d9\u05e9\u05d9\u05e8\u05d9\u05dd \u05e9\u05dc \u05e1\u05d9\u05de\u05e4\u05dc\u05e7\u05e1\u05d9\u05dd.\",\n  \"Replace\"\n);\n\n// 5) Third body paragraph\nitems[4].insertText(\n  \"\u05d4\u05de\u05d7\u05e7\u05e8 \u05d1\u05d5\u05d7\u05df 957 \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05e2\u05dd \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05de- WordNet \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05de\u05d5\u05d3\u05dc \u05d2'\u05de\u05d4. \u05d4\u05db\u05d5\u05ea\u05d1\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e1\u05de\u05e0\u05d8\u05d9\u05d9\u05dd high-level \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e0\u05d5\u05d8\u05e8\u05d9\u05dd \u05d5\u05de\u05e0\u05d5\u05d4\u05dc\u05d9\u05dd \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d3\u05d9\u05d3\u05d4 \u05d5\u05e2\u05e8\u05d9\u05db\u05d4 \u05d9\u05e9\u05d9\u05e8\u05d4 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05d9\u05dd \u05d4\u05e4\u05e0\u05d9\u05de\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4-LLMs. \u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05d2\u05dc\u05d5\u05ea \u05de\u05d1\u05e0\u05d4 \u05e4\u05e9\u05d5\u05d8 \u05e9\u05d1\u05d5 \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d9\u05dd \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05d2\u05d9\u05d0\u05d5\u05de\u05d8\u05e8\u05d9\u05ea \u05db\u05e1\u05d9\u05de\u05e4\u05dc\u05e7\u05e1\u05d9\u05dd \u05d5\u05de\u05d5\u05e9\u05d2\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05de\u05e7\u05d5\u05d3\u05d3\u05d9\u05dd \u05db\u05d0\u05d5\u05e8\u05ea\u05d5\u05d2\u05d5\u05e0\u05dc\u05d9\u05d5\u05ea.\",\n  \"Replace\"\n);\n\n// 6) Remove the paragraph that discussed autoregressive inference speed-up\n//    (no replacement text in the new version of the review).\nitems[5].delete();\n\n// 7) Update the arxiv link (now at index 5 after the delete above).\nitems[6].insertText(\"https://arxiv.org/pdf/2406.01506\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Apply the Review_221.docx edit via the Word COM object model:\n#  P1: date 12.06.24 -> 11.06.24\n#  P2: title replaced (new paper title)\n#  P3: body paragraph replaced\n#  P4: body paragraph replaced\n#  P5: body paragraph replaced\n#  P6: body paragraph removed entirely (old \"\u05d1\u05d3\u05f4\u05db \u05d4\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1...\" paragraph)\n#  P7 (becomes P6 after removal): URL replaced with the new arxiv link\n\n$d = $word.ActiveDocument\n\n# 1) Date line\n$d.Paragraphs(1).Range.Text = \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 11.06.24:\u26a1\ufe0f\ud83d\ude80\"\n\n# 2) Paper title\n$d.Paragraphs(2).Range.Text = \"The Geometry of Categorical and Hierarchical Concepts in Large Language Models\"\n\n# 3) First body paragraph\n$d.Paragraphs(3).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 \u05db\u05d9\u05e6\u05d3 \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05d5\u05de\u05d5\u05e9\u05d2\u05d9\u05dd \u05de\u05e7\u05d5\u05d3\u05d3\u05d9\u05dd \u05d1\u05de\u05e8\u05d7\u05d1\u05d9 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 (embeddings)\u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05dc \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d4. \u05d4\u05db\u05d5\u05ea\u05d1\u05d9\u05dd \u05d7\u05d5\u05e7\u05e8\u05d9\u05dd 2 \u05e9\u05d0\u05dc\u05d5\u05ea \u05de\u05e8\u05db\u05d6\u05d9\u05d5\u05ea: \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d9\u05dd \u05d5\u05d4\u05e7\u05d9\u05d3\u05d5\u05d3 \u05e9\u05dc \u05d9\u05d7\u05e1\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05d1\u05d9\u05df \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd.\"\n\n# 4) Second body paragraph\n$d.Paragraphs(4).Range.Text = \"\u05d4\u05dd \u05de\u05e8\u05d7\u05d9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05d4\u05e1\u05ea\u05db\u05dc\u05d5\u05ea \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d4\u05e8\u05d2\u05d9\u05dc\u05d4 \u05e2\u05dc \u05d4\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05d4\u05e8\u05d0\u05d5\u05ea \u05e9\u05d4\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d9\u05dd \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05db\u05e1\u05d9\u05de\u05e4\u05dc\u05e7\u05e1\u05d9\u05dd, \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05d4\u05dd \u05d0\u05d5\u05e8\u05ea\u05d5\u05d2\u05d5\u05e0\u05dc\u05d9\u05d9\u05dd, \u05d5\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05db\u05e4\u05d5\u05dc\u05d9\u05d8\u05d5\u05e4\u05d9\u05dd \u05e9\u05e0\u05d1\u05e0\u05d9\u05dd \u05de\u05e1\u05db\u05d5\u05de\u05d9\u05dd \u05d9\u05e9\u05d9\u05e8\u05d9\u05dd \u05e9\u05dc \u05e1\u05d9\u05de\u05e4\u05dc\u05e7\u05e1\u05d9\u05dd.\"\n\n# 5) Third body paragraph\n$d.Paragraphs(5).Range.Text = \"\u05d4\u05de\u05d7\u05e7\u05e8 \u05d1\u05d5\u05d7\u05df 957 \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05e2\u05dd \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05de- WordNet \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05de\u05d5\u05d3\u05dc \u05d2'\u05de\u05d4. \u05d4\u05db\u05d5\u05ea\u05d1\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e1\u05de\u05e0\u05d8\u05d9\u05d9\u05dd high-level \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e0\u05d5\u05d8\u05e8\u05d9\u05dd \u05d5\u05de\u05e0\u05d5\u05d4\u05dc\u05d9\u05dd \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d3\u05d9\u05d3\u05d4 \u05d5\u05e2\u05e8\u05d9\u05db\u05d4 \u05d9\u05e9\u05d9\u05e8\u05d4 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05d9\u05dd \u05d4\u05e4\u05e0\u05d9\u05de\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4-LLMs. \u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05d2\u05dc\u05d5\u05ea \u05de\u05d1\u05e0\u05d4 \u05e4\u05e9\u05d5\u05d8 \u05e9\u05d1\u05d5 \u05e7\u05d5\u05e0\u05e1\u05e4\u05d8\u05d9\u05dd \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d9\u05dd \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05d2\u05d9\u05d0\u05d5\u05de\u05d8\u05e8\u05d9\u05ea \u05db\u05e1\u05d9\u05de\u05e4\u05dc\u05e7\u05e1\u05d9\u05dd \u05d5\u05de\u05d5\u05e9\u05d2\u05d9\u05dd \u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d9\u05dd \u05de\u05e7\u05d5\u05d3\u05d3\u05d9\u05dd \u05db\u05d0\u05d5\u05e8\u05ea\u05d5\u05d2\u05d5\u05e0\u05dc\u05d9\u05d5\u05ea.\"\n\n# 6) Remove the paragraph that discussed autoregressive inference speed-up\n#    (no replacement text in the new version of the review).\n$d.Paragraphs(6).Range.Delete()\n\n# 7) Update the arxiv link (now paragraph 6 after the delete above).\n$d.Paragraphs(6).Range.Text = \"https://arxiv.org/pdf/2406.01506\"\n"}
